# Insert two brand-new weekly price rows ("Choclo" / "Choclero" and
# "Choclo" / "Dulce o Americano") at the top of the Vega Modelo de Temuco
# data block (rows 386-387), pushing all the existing rows down by two
# (previous A1:R475 used range becomes A1:R477).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 386:475 down to 388:477, inserting two blank rows.
$ws.Range("A386:R387").EntireRow.Insert()

# New row 386 - Choclo / Choclero / Primera
$ws.Range("A386").Value = 10
$ws.Range("B386").Value = "Vega Modelo de Temuco"
$ws.Range("C386").Value = "La Araucanía"
$ws.Range("D386").Value = 44889
$ws.Range("E386").Value = 9
$ws.Range("F386").Value = 100112024
$ws.Range("G386").Value = "Choclo"
$ws.Range("H386").Value = "Choclero"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 155
$ws.Range("K386").Value = 28000
$ws.Range("L386").Value = 28000
$ws.Range("M386").Value = 28000
$ws.Range("N386").Value = "$/malla 50 unidades"
$ws.Range("O386").Value = "Región de Arica y Parinacota"
$ws.Range("P386").Value = 560
$ws.Range("Q386").Value = 50
$ws.Range("R386").Value = "Hortaliza"

# New row 387 - Choclo / Dulce o Americano / Primera
$ws.Range("A387").Value = 10
$ws.Range("B387").Value = "Vega Modelo de Temuco"
$ws.Range("C387").Value = "La Araucanía"
$ws.Range("D387").Value = 44889
$ws.Range("E387").Value = 9
$ws.Range("F387").Value = 100112024
$ws.Range("G387").Value = "Choclo"
$ws.Range("H387").Value = "Dulce o Americano"
$ws.Range("I387").Value = "Primera"
$ws.Range("J387").Value = 250
$ws.Range("K387").Value = 28000
$ws.Range("L387").Value = 28000
$ws.Range("M387").Value = 28000
$ws.Range("N387").Value = "$/malla 70 unidades"
$ws.Range("O387").Value = "Región de Arica y Parinacota"
$ws.Range("P387").Value = 400
$ws.Range("Q387").Value = 70
$ws.Range("R387").Value = "Hortaliza"
